{"js": "// Update the worksheet date and every two-digit multiplication prompt.\n// Each entry is [oldText, newText]; all values in this document are\n// unique, so an exact, case-sensitive whole-match search finds exactly\n// one hit per entry.\nconst replacements = [\n  [\"2025-01-22 Wednesday\", \"2025-01-23 Thursday\"],\n  [\"81\u00d790=\", \"42\u00d792=\"],\n  [\"90\u00d712=\", \"23\u00d712=\"],\n  [\"87\u00d793=\", \"28\u00d768=\"],\n  [\"98\u00d783=\", \"20\u00d733=\"],\n  [\"21\u00d792=\", \"99\u00d731=\"],\n  [\"15\u00d746=\", \"18\u00d775=\"],\n  [\"83\u00d763=\", \"61\u00d757=\"],\n  [\"50\u00d732=\", \"65\u00d785=\"],\n  [\"65\u00d761=\", \"48\u00d783=\"],\n  [\"95\u00d743=\", \"78\u00d797=\"],\n  [\"77\u00d735=\", \"77\u00d749=\"],\n  [\"46\u00d711=\", \"74\u00d765=\"],\n  [\"36\u00d751=\", \"64\u00d757=\"],\n  [\"37\u00d766=\", \"77\u00d776=\"],\n  [\"94\u00d729=\", \"67\u00d761=\"],\n  [\"67\u00d715=\", \"17\u00d719=\"],\n  [\"89\u00d762=\", \"92\u00d720=\"],\n  [\"63\u00d754=\", \"59\u00d752=\"],\n  [\"59\u00d720=\", \"25\u00d713=\"],\n  [\"87\u00d740=\", \"22\u00d724=\"],\n  [\"32\u00d787=\", \"92\u00d725=\"],\n  [\"60\u00d734=\", \"54\u00d772=\"],\n  [\"69\u00d794=\", \"16\u00d780=\"],\n  [\"78\u00d719=\", \"20\u00d782=\"],\n  [\"84\u00d783=\", \"65\u00d734=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`No match found for \"${oldText}\"`);\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date and every two-digit multiplication prompt.\n# Each tuple is (oldText, newText); every value in this document is\n# unique, so a literal, whole-document Find/Replace locates exactly one\n# occurrence per pair.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @('2025-01-22 Wednesday', '2025-01-23 Thursday'),\n    @('81\u00d790=', '42\u00d792='),\n    @('90\u00d712=', '23\u00d712='),\n    @('87\u00d793=', '28\u00d768='),\n    @('98\u00d783=', '20\u00d733='),\n    @('21\u00d792=', '99\u00d731='),\n    @('15\u00d746=', '18\u00d775='),\n    @('83\u00d763=', '61\u00d757='),\n    @('50\u00d732=', '65\u00d785='),\n    @('65\u00d761=', '48\u00d783='),\n    @('95\u00d743=', '78\u00d797='),\n    @('77\u00d735=', '77\u00d749='),\n    @('46\u00d711=', '74\u00d765='),\n    @('36\u00d751=', '64\u00d757='),\n    @('37\u00d766=', '77\u00d776='),\n    @('94\u00d729=', '67\u00d761='),\n    @('67\u00d715=', '17\u00d719='),\n    @('89\u00d762=', '92\u00d720='),\n    @('63\u00d754=', '59\u00d752='),\n    @('59\u00d720=', '25\u00d713='),\n    @('87\u00d740=', '22\u00d724='),\n    @('32\u00d787=', '92\u00d725='),\n    @('60\u00d734=', '54\u00d772='),\n    @('69\u00d794=', '16\u00d780='),\n    @('78\u00d719=', '20\u00d782='),\n    @('84\u00d783=', '65\u00d734=')\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
